$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.11017566666667
$ws.Range("H2").Value = 39.330527
$ws.Range("I2").Value = 0.1657114824704502
$ws.Range("J2").Value = 0.1657114824704501
$ws.Range("M2").Value = 5.916202333333334
$ws.Range("N2").Value = 17.748607
$ws.Range("O2").Value = 0.3515586392055965
$ws.Range("P2").Value = 0.3515586392055965
$ws.Range("Q2").Value = 77.56245186954324
$ws.Range("R2").Value = 698.062066825889
$ws.Range("S2").Value = 0.05825730327805353
$ws.Range("T2").Value = 0.05825730327805352

$ws.Range("G3").Value = 13.11017566666667
$ws.Range("H3").Value = 39.330527
$ws.Range("I3").Value = 0.1657114824704502
$ws.Range("J3").Value = 0.1657114824704501
$ws.Range("O3").Value = 0.6159539016771971
$ws.Range("P3").Value = 0.6159539016771971
$ws.Range("Q3").Value = 135.8945266162426
$ws.Range("R3").Value = 1223.050739546183
$ws.Range("S3").Value = 0.1020706341803862
$ws.Range("T3").Value = 0.1020706341803862

$ws.Range("G4").Value = 13.11017566666667
$ws.Range("H4").Value = 39.330527
$ws.Range("I4").Value = 0.1657114824704502
$ws.Range("J4").Value = 0.1657114824704501
$ws.Range("M4").Value = 0.5467150000000001
$ws.Range("O4").Value = 0.03248745911720639
$ws.Range("P4").Value = 0.03248745911720639
$ws.Range("Q4").Value = 7.167529689601668
$ws.Range("R4").Value = 64.50776720641501
$ws.Range("S4").Value = 0.005383545012010414
$ws.Range("T4").Value = 0.005383545012010413

$ws.Range("G5").Value = 51.42568199999999
$ws.Range("I5").Value = 0.6500161567583834
$ws.Range("J5").Value = 0.6500161567583833
$ws.Range("M5").Value = 5.916202333333334
$ws.Range("N5").Value = 17.748607
$ws.Range("O5").Value = 0.3515586392055965
$ws.Range("P5").Value = 0.3515586392055965
$ws.Range("Q5").Value = 304.244739841658
$ws.Range("R5").Value = 2738.202658574922
$ws.Range("S5").Value = 0.228518795531629
$ws.Range("T5").Value = 0.228518795531629

$ws.Range("G6").Value = 51.42568199999999
$ws.Range("I6").Value = 0.6500161567583834
$ws.Range("J6").Value = 0.6500161567583833
$ws.Range("O6").Value = 0.6159539016771971
$ws.Range("P6").Value = 0.6159539016771971
$ws.Range("Q6").Value = 533.0568322647259
$ws.Range("R6").Value = 4797.511490382533
$ws.Range("S6").Value = 0.4003799879085428
$ws.Range("T6").Value = 0.4003799879085428

$ws.Range("G7").Value = 51.42568199999999
$ws.Range("I7").Value = 0.6500161567583834
$ws.Range("J7").Value = 0.6500161567583833
$ws.Range("M7").Value = 0.5467150000000001
$ws.Range("O7").Value = 0.03248745911720639
$ws.Range("P7").Value = 0.03248745911720639
$ws.Range("S7").Value = 0.0211173733182116
$ws.Range("T7").Value = 0.0211173733182116

$ws.Range("I8").Value = 0.1842723607711665
$ws.Range("J8").Value = 0.1842723607711665
$ws.Range("M8").Value = 5.916202333333334
$ws.Range("N8").Value = 17.748607
$ws.Range("O8").Value = 0.3515586392055965
$ws.Range("P8").Value = 0.3515586392055965
$ws.Range("Q8").Value = 86.25000452668921
$ws.Range("R8").Value = 776.250040740203
$ws.Range("S8").Value = 0.06478254039591405
$ws.Range("T8").Value = 0.06478254039591404

$ws.Range("I9").Value = 0.1842723607711665
$ws.Range("J9").Value = 0.1842723607711665
$ws.Range("O9").Value = 0.6159539016771971
$ws.Range("P9").Value = 0.6159539016771971
$ws.Range("S9").Value = 0.1135032795882681
$ws.Range("T9").Value = 0.1135032795882681

$ws.Range("I10").Value = 0.1842723607711665
$ws.Range("J10").Value = 0.1842723607711665
$ws.Range("M10").Value = 0.5467150000000001
$ws.Range("O10").Value = 0.03248745911720639
$ws.Range("P10").Value = 0.03248745911720639
$ws.Range("Q10").Value = 7.970344583911666
$ws.Range("R10").Value = 71.733101255205
$ws.Range("S10").Value = 0.00598654078698438
$ws.Range("T10").Value = 0.005986540786984378
